# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

# Row number -> [old value, new value] for column F
$updates = @{
    3  = 2188
    4  = 77
    5  = 12916
    8  = 509
    9  = 473
    10 = 1166
    11 = 968
    12 = 13722
    13 = 14206
    18 = 30
    22 = 1081
    26 = 5298
    28 = 285
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
